$wb = $excel.ActiveWorkbook

# ---- Section_A ----
$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("B2").Value = 'ELECTIVE_B1 [C204]'
$ws.Range("C2").Value = 'MA161 [C404]'
$ws.Range("D2").Value = 'ELECTIVE_B1 [C204]'
$ws.Range("E2").Value = 'Free'
$ws.Range("B3").Value = 'Free'
$ws.Range("C3").Value = 'DS161 [C302]'
$ws.Range("D3").Value = 'MA162 [C205]'
$ws.Range("E3").Value = 'MA162 [C205]'
$ws.Range("F3").Value = 'DS161 [C302]'
$ws.Range("C5").Value = 'EC161 (Lab) [L306]'
$ws.Range("D5").Value = 'Free'
$ws.Range("E5").Value = 'EC161 [C102]'
$ws.Range("F5").Value = 'Free'
$ws.Range("C6").Value = 'EC161 (Lab) [L306]'
$ws.Range("D6").Value = 'Free'
$ws.Range("B7").Value = 'MA161 [C404]'
$ws.Range("C7").Value = 'EC161 [C102]'
$ws.Range("D7").Value = 'Free'
$ws.Range("E7").Value = 'Free'

# ---- Section_B ----
$ws = $wb.Worksheets.Item("Section_B")
$ws.Range("B2").Value = 'ELECTIVE_B1 [C304]'
$ws.Range("D2").Value = 'ELECTIVE_B1 [C304]'
$ws.Range("E2").Value = 'DS161 [C004]'
$ws.Range("F2").Value = 'EC161 [C101]'
$ws.Range("D3").Value = 'MA161 [C202]'
$ws.Range("E3").Value = 'Free'
$ws.Range("F3").Value = 'MA162 [C101]'
$ws.Range("B5").Value = 'Free'
$ws.Range("C5").Value = 'MA162 [C101]'
$ws.Range("D5").Value = 'EC161 (Lab) [L408]'
$ws.Range("F5").Value = 'Free'
$ws.Range("D6").Value = 'EC161 (Lab) [L408]'
$ws.Range("F6").Value = 'Free'
$ws.Range("B7").Value = 'EC161 [C101]'
$ws.Range("C7").Value = 'DS161 [C004]'
$ws.Range("D7").Value = 'Free'
$ws.Range("F7").Value = 'MA161 [C202]'

# ---- Classroom_Utilization ----
$ws = $wb.Worksheets.Item("Classroom_Utilization")
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 0.6
$ws.Range("G5").Value = 7.5
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 1.2
$ws.Range("G6").Value = 15
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0.6
$ws.Range("G7").Value = 7.5
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0.6
$ws.Range("G16").Value = 7.5
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 0.6
$ws.Range("G22").Value = 7.5
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("D26").Value = 2.5
$ws.Range("E26").Value = 0.5
$ws.Range("G26").Value = 6.25
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("D32").Value = 3
$ws.Range("E32").Value = 0.6
$ws.Range("G32").Value = 7.5
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("D36").Value = 2.5
$ws.Range("E36").Value = 0.5
$ws.Range("G36").Value = 6.25

# ---- Classroom_Allocation ----
$ws = $wb.Worksheets.Item("Classroom_Allocation")
$ws.Range("G2").Value = 'C204'
$ws.Range("E3").Value = '15:30-17:00'
$ws.Range("F3").Value = 'MA161'
$ws.Range("G3").Value = 'C404'
$ws.Range("I3").Value = '''78'
$ws.Range("E4").Value = '09:00-10:30'
$ws.Range("G4").Value = 'C404'
$ws.Range("I4").Value = '''78'
$ws.Range("D5").Value = 'Tue'
$ws.Range("E5").Value = '10:30-12:00'
$ws.Range("F5").Value = 'DS161'
$ws.Range("G5").Value = 'C302'
$ws.Range("D6").Value = 'Tue'
$ws.Range("E6").Value = '13:00-14:30'
$ws.Range("F6").Value = 'EC161 (Lab)'
$ws.Range("G6").Value = 'L306'
$ws.Range("J6").Value = 'Computers'
$ws.Range("D7").Value = 'Tue'
$ws.Range("E7").Value = '14:30-15:30'
$ws.Range("G7").Value = 'L306'
$ws.Range("I7").Value = '''96'
$ws.Range("D8").Value = 'Tue'
$ws.Range("E8").Value = '15:30-17:00'
$ws.Range("F8").Value = 'EC161'
$ws.Range("G8").Value = 'C102'
$ws.Range("I8").Value = '''96'
$ws.Range("J8").Value = 'Projector'
$ws.Range("E9").Value = '09:00-10:30'
$ws.Range("F9").Value = 'ELECTIVE_B1'
$ws.Range("G9").Value = 'C204'
$ws.Range("D10").Value = 'Wed'
$ws.Range("E10").Value = '10:30-12:00'
$ws.Range("F10").Value = 'MA162'
$ws.Range("G10").Value = 'C205'
$ws.Range("E11").Value = '10:30-12:00'
$ws.Range("F11").Value = 'MA162'
$ws.Range("G11").Value = 'C205'
$ws.Range("D12").Value = 'Thu'
$ws.Range("E12").Value = '13:00-14:30'
$ws.Range("G12").Value = 'C102'
$ws.Range("E13").Value = '10:30-12:00'
$ws.Range("F13").Value = 'DS161'
$ws.Range("G13").Value = 'C302'
$ws.Range("G14").Value = 'C304'
$ws.Range("E15").Value = '15:30-17:00'
$ws.Range("F15").Value = 'EC161'
$ws.Range("G15").Value = 'C101'
$ws.Range("I15").Value = '''96'
$ws.Range("F16").Value = 'MA162'
$ws.Range("G16").Value = 'C101'
$ws.Range("F17").Value = 'DS161'
$ws.Range("G17").Value = 'C004'
$ws.Range("H17").Value = 'Auditorium'
$ws.Range("I17").Value = '''240'
$ws.Range("G18").Value = 'C304'
$ws.Range("F19").Value = 'MA161'
$ws.Range("G19").Value = 'C202'
$ws.Range("I19").Value = '''96'
$ws.Range("E20").Value = '13:00-14:30'
$ws.Range("F20").Value = 'EC161 (Lab)'
$ws.Range("G20").Value = 'L408'
$ws.Range("H20").Value = 'classroom without projector'
$ws.Range("I20").Value = '''78'
$ws.Range("J20").Value = 'Computers'
$ws.Range("D21").Value = 'Wed'
$ws.Range("E21").Value = '14:30-15:30'
$ws.Range("F21").Value = 'EC161 (Lab)'
$ws.Range("G21").Value = 'L408'
$ws.Range("H21").Value = 'classroom without projector'
$ws.Range("I21").Value = '''78'
$ws.Range("J21").Value = 'Computers'
$ws.Range("D22").Value = 'Thu'
$ws.Range("G22").Value = 'C004'
$ws.Range("H22").Value = 'Auditorium'
$ws.Range("I22").Value = '''240'
$ws.Range("E23").Value = '09:00-10:30'
$ws.Range("G23").Value = 'C101'
$ws.Range("E24").Value = '10:30-12:00'
$ws.Range("F24").Value = 'MA162'
$ws.Range("G24").Value = 'C101'
$ws.Range("I24").Value = '''96'
$ws.Range("J24").Value = 'Projector'
$ws.Range("E25").Value = '15:30-17:00'
$ws.Range("F25").Value = 'MA161'
$ws.Range("G25").Value = 'C202'
$ws.Range("I25").Value = '''96'
$ws.Range("J25").Value = 'Projector'

